$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cpu_windows log entries appended after the existing data (rows 2-81).
# Each tuple is (date_hour label, cpu_used value); the label cycles through the
# same "day_hour" vocabulary already present in column A, so Excel reuses the
# matching shared-string entry whenever the label already exists.
$logEntries = @(
    ,@("6_12", 97.40000000000001)
    ,@("6_18", 98.09)
    ,@("7_00", 97.06999999999999)
    ,@("7_06", 98.59)
    ,@("7_12", 9.01)
    ,@("7_18", 13.56)
    ,@("8_00", 10.28)
    ,@("8_06", 16.54)
    ,@("8_12", 25.98)
    ,@("8_18", 21.21)
    ,@("9_00", 20.45)
    ,@("9_06", 19.79)
    ,@("9_12", 26.43)
    ,@("9_18", 25.83)
    ,@("10_00", 16.95)
    ,@("10_06", 7.61)
    ,@("10_12", 41.74)
    ,@("10_18", 26.73)
    ,@("11_00", 35.94)
    ,@("11_06", 65.8)
    ,@("11_12", 40.08)
    ,@("11_18", 46.8)
    ,@("12_00", 25.16)
    ,@("12_06", 23.83)
    ,@("12_12", 82.70999999999999)
    ,@("12_18", 84.90000000000001)
    ,@("13_00", 27.37)
    ,@("13_06", 28.47)
    ,@("13_12", 25.64)
    ,@("13_18", 32.86)
    ,@("14_00", 34.45)
    ,@("14_06", 41.28)
    ,@("14_12", 92.56)
    ,@("14_18", 53.51)
    ,@("15_00", 45.13)
    ,@("15_06", 44.37)
    ,@("15_12", 41.91)
    ,@("15_18", 37.26)
    ,@("16_00", 32.85)
    ,@("16_06", 50.17)
    ,@("16_12", 32.62)
    ,@("16_18", 25.67)
    ,@("17_00", 22.37)
    ,@("17_06", 22.42)
    ,@("17_12", 15.34)
    ,@("17_18", 12.12)
    ,@("18_00", 11.39)
    ,@("18_06", 22.4)
    ,@("18_12", 19.3)
    ,@("18_18", 16.65)
    ,@("19_00", 18.44)
    ,@("19_06", 15.51)
    ,@("19_12", 16.97)
    ,@("19_18", 25.59)
    ,@("20_00", 20.24)
    ,@("20_06", 15.62)
)

$startRow = $ws.UsedRange.Rows.Count + 1
$r = $startRow
foreach ($entry in $logEntries) {
    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $r++
}
